# ---------------------------------------------------------------------------
# Commit: "agregado de conf. git"
#
# The only change recorded by this commit's OOXML diff is inside
# ppt/changesInfos/changesInfo1.xml: two already-present <pc:grpChg> history
# records for the slide-master layout's group shape (shape id 111, under
# sldMasterMk cId="702745799" / sldLayoutMk sldId="2147483660") are
# reordered -- the "add mod" record for creationId {1CC7FA66-...} moves to
# sit right before the "del" record for creationId {87376B41-...} (same
# two <pc:grpChg> elements, same attributes/timestamps/values, just swapped
# relative order).
#
# No shape is added, deleted, moved, resized, restyled or retexted; no
# slide, layout, master, theme or media content changes anywhere in the
# deck -- the live slideLayout2.xml already reflects the {1CC7FA66-...}
# group (id 111) as the current shape, exactly as before.
#
# ppt/changesInfos/changesInfo1.xml is PowerPoint's internal
# coauthoring/version-history changelog. It isn't part of the
# Presentation/Slides/Shapes automation object model -- there is no
# property or method on $ppt.ActivePresentation (or anything reachable
# from it) that reads or writes that log, so a COM/VBA macro has no way to
# poke its entries, and no legitimate action here would do so as a side
# effect either. Since nothing in the visible/editable content model
# changed, the correct macro is a no-op against that model: touching any
# real shape would introduce content differences the diff never asked for.
$p = $ppt.ActivePresentation
